# Append a new case row (row 6) to the "ランサーズ" sheet and refresh the
# "取得日時" (fetched-at) timestamp on every existing data row, per commit:
#   "Append: 2026-02-14 12:42 JST"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-14 12:42:13"

# Drop every hyperlink up front -- this engine's Hyperlinks.Delete() always
# clears the whole worksheet collection, and Rows.Insert() does not relocate
# existing hyperlink anchors, so it is simplest/most reliable to rebuild the
# hyperlinks once all the row data is in its final place.
$ws.Hyperlinks.Delete()

# Insert a brand-new row above the old row 6, pushing the two Github/Vercel
# rows down to rows 7 and 8 (matching the diff's row shift).
$ws.Rows(6).Insert()

# --- New row 6: the C2C skill-marketplace listing ---
$ws.Range("B6").Value = "【急募】C2Cスキル売買プラットフォーム構築の依頼"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5491832"
$ws.Range("G6").Value = 25

# Column D got a little wider to fit the bigger price range text.
# (ColumnWidth uses Excel's character-width units, which store as
# value + 6/7 in the saved file, so back that padding out here.)
$ws.Columns("D").ColumnWidth = 32 - 6/7

# Refresh the "取得日時" timestamp across every data row (2-8), including the
# two rows that just shifted down and the freshly inserted row 6.
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Rebuild hyperlinks top to bottom so relationship ids line up the same way
# Excel would renumber them (F2..F5 unchanged content, F6 new, F7/F8 are the
# shifted-down rows).
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5491704")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5491672")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5491578")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5491569")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5491832")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5491736")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5491643")

# Hyperlinks.Add() nudges the cell onto a freshly-minted duplicate of the
# built-in "Hyperlink" style; re-apply the named style so every link cell
# collapses back onto the single shared Hyperlink style record.
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
}
